# Updates the "cryptos" price/volume snapshot in columns D (Price) and E
# (Volume 1h) for rows 2-51, matching the refreshed scrape values from the
# "Updated cryptos list ... with GitHub Actions" run.
#
# Some new Price values (column D) are plain decimal numbers (e.g. "1.003").
# Assigning such a string straight to .Value lets Excel auto-coerce it to a
# numeric cell, which would change both the stored type and drop the leading/
# trailing formatting of values like "0.000008723". To keep these as literal
# text (as they were before the edit) we briefly force a text number format,
# write the value, then restore the cell style so no stray formatting is left
# behind. Values that already can't parse as a plain number (e.g.
# "27.883.33", which has two dots) are assigned directly since Excel keeps
# those as text on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.883.33'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '1.908.14'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5017'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3820'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9099'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.46%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '1.935.91'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07669'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008723'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = '27.921.94'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.171'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.587'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("E25").Value = '  -2.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.224'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.916'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09006'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.215'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.230'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7633'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.655'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02062'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("E36").Value = '  -4.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5583'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.092'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.022'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05251'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.938'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.498'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1512'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4839'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.629'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9028'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.26%  '
